# CIR-1787 - Updated XLS/CVS service guide for story
#
# The "END_DATE_AFTER_GROUP_END_DATE" error code (and its message) is no
# longer a valid error for the Full Return, so remove its row from the
# "Full Return Errors" sheet. Deleting the whole row shifts every
# subsequent row up by one, which is exactly what the published workbook
# shows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Full Return Errors")

$targetRow = $null
$lastRow = $ws.Cells(1, 1).End(4).Row  # xlDown = 4 -> find used range bottom via column A

for ($r = 1; $r -le $lastRow; $r++) {
    if ($ws.Cells.Item($r, 1).Value() -eq "END_DATE_AFTER_GROUP_END_DATE") {
        $targetRow = $r
        break
    }
}

if ($targetRow -ne $null) {
    $ws.Rows.Item($targetRow).Delete()
}

$ws.Range("A1").Select()
